$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.294.86'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.588.58'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  -0.80%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.10'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.503'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('E8').Value = '  +0.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0610'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.34'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '1.812.81'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '1.613.28'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.35'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').Value = '26.301.95'
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.43'
$ws.Range('E19').Value = '  +5.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '210.83'
$ws.Range('E20').Value = '  +2.91%  '
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.27'
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.94'
$ws.Range('E23').Value = '  +1.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.14'
$ws.Range('E24').Value = '  -2.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.57'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.03'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.23'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.99'
$ws.Range('E33').Value = '  +1.91%  '
$ws.Range('D34').Value = '1.318.35'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.603'
$ws.Range('E36').Value = '  +2.03%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0167'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('E39').Value = '  -13.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.808'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  +4.82%  '
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.14'
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.01'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '1.724.78'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '87.61'
$ws.Range('E47').Value = '  -0.98%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.48'
$ws.Range('E49').Value = '  -4.91%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0505'
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0973'
$ws.Range('E51').Value = '  -3.71%  '
